# "Generate Report for Handoff" — a fresh localization-status report run
# stamped a new "Latest Handoff Datetime" for the files that were handed
# off in this batch (zh-cn: 2016-03-10 12:36:22, de-de: 2016-03-10 12:36:27).
#
# Rows 7, 10-16 on each language sheet are all part of this handoff batch.
# Previously rows 10 & 11 ("05883ead..." / "0f2fb72c...") carried their own
# slightly-earlier timestamp (12:36:04 / 12:36:08); the new report run gives
# every row in the batch the same, latest timestamp.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$zhTimestamp = "2016-03-10 12:36:22"
$deTimestamp = "2016-03-10 12:36:27"

$batchRows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $batchRows) {
    $zh.Range("D$r").Value = $zhTimestamp
    $de.Range("D$r").Value = $deTimestamp
}
